$d = $word.ActiveDocument

$replacements = @(
    @{old="427×9=3843"; new="303×5=1515"},
    @{old="882×4=3528"; new="208×9=1872"},
    @{old="835×8=6680"; new="865×5=4325"},
    @{old="637×8=5096"; new="145×2=290"},
    @{old="491×9=4419"; new="641×2=1282"},
    @{old="620×7=4340"; new="304×8=2432"},
    @{old="142×6=852"; new="864×7=6048"},
    @{old="964×2=1928"; new="398×4=1592"},
    @{old="528×5=2640"; new="456×6=2736"},
    @{old="677×2=1354"; new="744×8=5952"},
    @{old="847×3=2541"; new="402×3=1206"},
    @{old="939×3=2817"; new="365×3=1095"},
    @{old="812×2=1624"; new="857×4=3428"},
    @{old="286×8=2288"; new="218×2=436"},
    @{old="283×6=1698"; new="745×3=2235"},
    @{old="690×6=4140"; new="713×2=1426"},
    @{old="996×4=3984"; new="456×6=2736"},
    @{old="911×6=5466"; new="852×6=5112"},
    @{old="771×8=6168"; new="422×8=3376"},
    @{old="395×9=3555"; new="664×7=4648"},
    @{old="983×9=8847"; new="736×5=3680"},
    @{old="202×7=1414"; new="927×4=3708"},
    @{old="240×9=2160"; new="127×6=762"},
    @{old="526×4=2104"; new="992×3=2976"},
    @{old="700×3=2100"; new="843×2=1686"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
